$d = $word.ActiveDocument
$s = $d.Styles.Add("Encabezado", 1)
try { $d.Save() } catch { Write-Host "SAVE ERR: $_" }
$s2 = $d.Styles.Item("Encabezado")
$ts = $s2.ParagraphFormat.TabStops
Write-Host "count before: $($ts.Count)"
$ts.ClearAll()
$ts.Add(4252, 1, 0)
Write-Host "count after: $($ts.Count)"
